$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the "INTER_REC" main/sub scenario rows (previously rows 16 & 17),
# shifting all the rows below them up by two.
$ws.Rows("16:17").Delete()

# Excel leaves the selection on the cell the user was last working on.
$ws.Range("D13").Select()
